$p = $ppt.ActivePresentation

# --- Slide 1 ("Parallel hardware architectures" title slide) ---
# Subtitle placeholder: "Fall 2015" -> "Fall " + "2016" (course refreshed to Fall 2016)
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2)
$subtitleTextRange = $subtitle.TextFrame.TextRange
$fallParagraph = $subtitleTextRange.Paragraphs(2, 1)
$yearChars = $fallParagraph.Characters(6, 4)
$yearChars.Text = "2016"

# --- Slide 2 (Chapter reference slide) ---
# Merge the three separate "Chapter " / "6 " / "from Computer Organization and Design"
# runs back into a single run with the same text/hyperlink.
$s2 = $p.Slides.Item(2)
$content = $s2.Shapes.Item(1)
$contentTextRange = $content.TextFrame.TextRange
$chapterParagraph = $contentTextRange.Paragraphs(1, 1)
# Collapse to a placeholder run first so the re-assignment below rebuilds the
# paragraph as a single run (a same-text re-assignment alone is a no-op and
# keeps the original three runs).
$chapterParagraph.Text = "X"
$chapterParagraph = $contentTextRange.Paragraphs(1, 1)
$chapterParagraph.Text = "Chapter 6 from Computer Organization and Design"
